# Insert a new "Days per week" column into the telework case table (Table3),
# between "Reason" and "Short description".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the existing table and its style before we touch anything, then take
# it out of "table" mode so we can freely restructure the underlying grid
# (ListObjects in this host don't re-sync column names/order after a raw
# column insert, so we rebuild the table once the grid is correct).
$tbl = $ws.ListObjects.Item(1)
$tableName = $tbl.Name
$tableStyle = $tbl.TableStyle.Name
$tbl.Unlist()

# Match the new column's width to the "Reason" column before inserting, so
# the inserted column does not inherit any bestFit/autosized behaviour.
$reasonWidth = $ws.Columns("D").ColumnWidth

# Insert a new blank column E; everything from old E onward (Short
# description, State, Assigned to) shifts one column to the right.
$ws.Columns("E").Insert()

# Header for the new column.
$ws.Range("E1").Value = "Days per week"

# Rows 2 and 3 had their Arrangement/Reason values swapped as part of this
# edit (independent of the new column).
$c2 = $ws.Range("C2").Value()
$d2 = $ws.Range("D2").Value()
$c3 = $ws.Range("C3").Value()
$d3 = $ws.Range("D3").Value()
$ws.Range("C2").Value = $c3
$ws.Range("D2").Value = $d3
$ws.Range("C3").Value = $c2
$ws.Range("D3").Value = $d2

# "Days per week" values captured from the source edit; every other data
# row is left blank in the new column.
$ws.Range("E2").Value = 1
$ws.Range("E9").Value = 3
$ws.Range("E15").Value = 5
$ws.Range("E16").Value = 2
$ws.Range("E17").Value = 1
$ws.Range("E18").Value = 3
$ws.Range("E19").Value = 1
$ws.Range("E20").Value = 2
$ws.Range("E24").Value = 1
$ws.Range("E25").Value = 3
$ws.Range("E30").Value = 2
$ws.Range("E31").Value = 5

# Column E takes a fixed width equal to the "Reason" column's width.
$ws.Columns("E").ColumnWidth = $reasonWidth

# Re-create the table over the now-correct A1:H40 range so ListObject
# column names/order are (re)derived from the header row.
$lastRow = $ws.Cells.Item(1, 1).End(-4121).Row
$newTbl = $ws.ListObjects.Add(1, $ws.Range("A1:H" + $lastRow), $null, 1)
$newTbl.Name = $tableName
$newTbl.TableStyle = $tableStyle

# Restore the user's selection roughly where it ended up after the edit.
$ws.Range("M22").Select()
